# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-column suffixes to the concrete
#    release tags "_FV2210" / "_FV2304".
# 2) Freeze the header row (row 1).
# 3) Turn the used range into a real Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -----------------------------------------------
$fields = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $fields[$i] + "_FV2210"
    $ws.Range($rightCols[$i] + "1").Value = $fields[$i] + "_FV2304"
}
# K1 stays "diff"

# --- 2) Freeze the top row --------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Convert the used range into an Excel Table ---------------------
$tableRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
